# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF holds the game "Date" as text in the format MM-DD-YYYY-YY
# (e.g. "12-24-2012-13"); it must instead read "2012-12-24" (YYYY-MM-DD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 contain the erroneous date strings (BF1 is the "Date" header).
$range = $ws.Range("BF2:BF31")

# Force text formatting first so Excel doesn't reinterpret the
# YYYY-MM-DD-looking string as a date serial value.
$range.NumberFormat = "@"
$range.Value = "2012-12-24"
